$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: predidx / pred_name updated (drop 'SurroundingEnvironment')
$ws.Range("D2").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E2").Value = "['Normal']"

# Row 56: predidx / pred_name updated (add 'HardwareFault')
$ws.Range("D56").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['HardwareFault']"

# Row 86: predidx / pred_name updated (drop 'ParamViolation')
$ws.Range("D86").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E86").Value = "['Normal']"
